$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-29 02:29:40"
$wsZh.Range("G2").Value = "2016-01-29 02:30:26"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-29 02:29:51"
$wsDe.Range("G2").Value = "2016-01-29 02:30:45"
